$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Steps"
$ws2 = $wb.Worksheets.Item(2)   # "Test Cases"

# --- Fix row 9 on "Test Steps": it used to be a stray "TS008 / Click on
# login button / btn_login / click" row; it becomes the first step
# (TS001 / Open the browser / openBrowser) of the TC002 block. ---
$ws1.Cells.Item(9, 2).Value = "TS001"
$ws1.Cells.Item(9, 3).Value = "Open the browser"
$ws1.Cells.Item(9, 4).Value = ""
$ws1.Cells.Item(9, 5).Value = "openBrowser"

# --- Bring formatting (borders) for the six new rows (10-15) in one
# shot, reusing the existing style used by the rest of the table so no
# new style entries are introduced. ---
$ws1.Range("A9:E9").Copy()
$ws1.Range("A10:E15").PasteSpecial(-4122)

# --- Populate the new rows: TC002's remaining test steps, mirroring
# TC001's steps (rows 3-8). ---
$ws1.Cells.Item(10, 1).Value = "TC002"
$ws1.Cells.Item(10, 2).Value = "TS002"
$ws1.Cells.Item(10, 3).Value = "Navigate to the url"
$ws1.Cells.Item(10, 4).Value = ""
$ws1.Cells.Item(10, 5).Value = "navigateUrl"

$ws1.Cells.Item(11, 1).Value = "TC002"
$ws1.Cells.Item(11, 2).Value = "TS003"
$ws1.Cells.Item(11, 3).Value = "Enter username"
$ws1.Cells.Item(11, 4).Value = "txt_name"
$ws1.Cells.Item(11, 5).Value = "input_uname"

$ws1.Cells.Item(12, 1).Value = "TC002"
$ws1.Cells.Item(12, 2).Value = "TS004"
$ws1.Cells.Item(12, 3).Value = "Enter password"
$ws1.Cells.Item(12, 4).Value = "txt_pass"
$ws1.Cells.Item(12, 5).Value = "input_upass"

$ws1.Cells.Item(13, 1).Value = "TC002"
$ws1.Cells.Item(13, 2).Value = "TS005"
$ws1.Cells.Item(13, 3).Value = "Click on login button"
$ws1.Cells.Item(13, 4).Value = "btn_login"
$ws1.Cells.Item(13, 5).Value = "click"

$ws1.Cells.Item(14, 1).Value = "TC002"
$ws1.Cells.Item(14, 2).Value = "TS006"
$ws1.Cells.Item(14, 3).Value = "Wait for sometime"
$ws1.Cells.Item(14, 4).Value = ""
$ws1.Cells.Item(14, 5).Value = "waitSometime"

$ws1.Cells.Item(15, 1).Value = "TC002"
$ws1.Cells.Item(15, 2).Value = "TS007"
$ws1.Cells.Item(15, 3).Value = "Close the browser"
$ws1.Cells.Item(15, 4).Value = ""
$ws1.Cells.Item(15, 5).Value = "closeBrowser"

# --- "Test Cases" sheet: TC002's Runmode flips from "No" to "Yes". ---
$ws2.Cells.Item(3, 3).Value = "Yes"

# --- Selection / active sheet: "Test Steps" is no longer the active
# tab; "Test Cases" becomes active with G3 selected. ---
$ws1.Range("A9").Select()
$ws2.Range("G3").Select()
